# Contest 1 MI vs RCB
# Update entered points for the first contest and fill in the upcoming
# contest match-ups, plus rename a couple of player header labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rename player headers in the merged row-11 cells ---
$ws.Range("Y11").Value = "Vjvignesh94"
$ws.Range("V11").Value = "Supriser Lee"
$ws.Range("AB11").Value = "Rag Nat"

# --- Contest 1 (row 13, "MI vs RCB") updated entered points ---
$ws.Range("E13").Value = 0
$ws.Range("H13").Value = 30
$ws.Range("K13").Value = 50
$ws.Range("N13").Value = 20
$ws.Range("Q13").Value = 80
$ws.Range("T13").Value = 60
$ws.Range("W13").Value = 70
$ws.Range("Z13").Value = 100
$ws.Range("AC13").Value = 40

# --- Fill in match-ups for upcoming contests (rows 14-27, column C) ---
$ws.Range("C14").Value = "CSK vs DC"
$ws.Range("C15").Value = "SRH vs KKR"
$ws.Range("C16").Value = "RR vs PBKS"
$ws.Range("C17").Value = "KKR vs MI"
$ws.Range("C18").Value = "SRH vs RCB"
$ws.Range("C19").Value = "RR vs DC"
$ws.Range("C20").Value = "PBKS vs CSK"
$ws.Range("C21").Value = "MI vs SRH"
$ws.Range("C22").Value = "RCB vs KKR"
$ws.Range("C23").Value = "DC vs PBKS"
$ws.Range("C24").Value = "CSK vs RR"
$ws.Range("C25").Value = "DC vs MI"
$ws.Range("C26").Value = "PBKS vs SRH"
$ws.Range("C27").Value = "KKR vs CSK"

$excel.CalculateFullRebuild()
